$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.091.33'
$ws.Range("E2").Value = '  -3.08%  '
$ws.Range("D3").Value = '2.520.41'
$ws.Range("E3").Value = '  -4.27%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'578.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").Value = "'168.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.98%  '
$ws.Range("D8").Value = "'0.520"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").Value = '2.520.73'
$ws.Range("E9").Value = '  -4.26%  '
$ws.Range("D10").Value = "'0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.82%  '
$ws.Range("E11").Value = '  -1.66%  '
$ws.Range("D12").Value = "'0.347"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.10%  '
$ws.Range("D13").Value = "'4.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.94%  '
$ws.Range("D14").Value = '2.969.73'
$ws.Range("E14").Value = '  -4.62%  '
$ws.Range("D15").Value = '69.793.33'
$ws.Range("E15").Value = '  -3.28%  '
$ws.Range("D16").Value = "'0.0000175"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.26%  '
$ws.Range("D17").Value = "'25.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.53%  '
$ws.Range("D18").Value = '2.520.94'
$ws.Range("E18").Value = '  -4.16%  '
$ws.Range("E19").Value = '  -0.60%  '
$ws.Range("D20").Value = "'11.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.65%  '
$ws.Range("D21").Value = "'350.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.32%  '
$ws.Range("D22").Value = "'3.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.96%  '
$ws.Range("D23").Value = "'1.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.04%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").Value = "'68.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.28%  '
$ws.Range("D26").Value = "'3.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.40%  '
$ws.Range("D27").Value = "'9.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.81%  '
$ws.Range("D28").Value = '2.647.45'
$ws.Range("E28").Value = '  -4.39%  '
$ws.Range("D29").Value = "'0.995"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("D30").Value = '0.0₃0911'
$ws.Range("E30").Value = '  -3.52%  '
$ws.Range("D31").Value = "'7.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("D32").Value = "'1.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.11%  '
$ws.Range("D33").Value = "'467.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.20%  '
$ws.Range("D34").Value = "'1.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.53%  '
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.12%  '
$ws.Range("D36").Value = "'0.119"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.68%  '
$ws.Range("D37").Value = "'152.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.98%  '
$ws.Range("D38").Value = "'19.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.75%  '
$ws.Range("D39").Value = "'18.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.13%  '
$ws.Range("D41").Value = "'4.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.79%  '
$ws.Range("D42").Value = "'0.320"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.35%  '
$ws.Range("D43").Value = "'1.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.46%  '
$ws.Range("D44").Value = "'1.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -13.57%  '
$ws.Range("D45").Value = "'2.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -9.24%  '
$ws.Range("D46").Value = "'38.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.02%  '
$ws.Range("D47").Value = "'143.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.57%  '
$ws.Range("D48").Value = "'0.532"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.84%  '
$ws.Range("D49").Value = "'3.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.50%  '
$ws.Range("E50").Value = '  -3.65%  '
$ws.Range("D51").Value = "'0.0735"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.20%  '
